# Update LR-pairs data (Ccl25-Ackr4) per Dr Hou's revised NATMI output (rows 2-13)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs", "Ccl25", "Ackr4", "ECs", 3, 1, 5.485012999999999, 16.455039, 0.2716272065325074, 0.2716272065325074, 3, 1, 0.2315733333333333, 0.69472, 0.2371078251520917, 0.2371078251520917, 1.270182743786666, 11.43164469408, 0.06440493619306087, 0.06440493619306087),
    @("ECs", "Ccl25", "Ackr4", "FAPs", 3, 1, 5.485012999999999, 16.455039, 0.2716272065325074, 0.2716272065325074, 3, 1, 0.7200953333333334, 2.160286, 0.7373052671097876, 0.7373052671097876, 3.949732264572667, 35.547590381154, 0.2002721700667358, 0.2002721700667358),
    @("ECs", "Ccl25", "Ackr4", "sCs", 3, 1, 5.485012999999999, 16.455039, 0.2716272065325074, 0.2716272065325074, 1, 0.3333333333333333, 0.02498966666666666, 0.07496899999999999, 0.02558690773812063, 0.02558690773812063, 0.1370686465323333, 1.233617818791, 0.006950100272710703, 0.006950100272710704),
    @("FAPs", "Ccl25", "Ackr4", "ECs", 3, 1, 6.542812333333333, 19.628437, 0.3240112351547335, 0.3240112351547334, 3, 1, 0.2315733333333333, 0.69472, 0.2371078251520917, 0.2371078251520917, 1.515140861404444, 13.63626775264, 0.07682559929238182, 0.07682559929238181),
    @("FAPs", "Ccl25", "Ackr4", "FAPs", 3, 1, 6.542812333333333, 19.628437, 0.3240112351547335, 0.3240112351547334, 3, 1, 0.7200953333333334, 2.160286, 0.7373052671097876, 0.7373052671097876, 4.711448628109111, 42.403037652982, 0.238895190282333, 0.2388951902823329),
    @("FAPs", "Ccl25", "Ackr4", "sCs", 3, 1, 6.542812333333333, 19.628437, 0.3240112351547335, 0.3240112351547334, 1, 0.3333333333333333, 0.02498966666666666, 0.07496899999999999, 0.02558690773812063, 0.02558690773812063, 0.1635026992725555, 1.471524293453, 0.008290445580018671, 0.008290445580018671),
    @("M2", "Ccl25", "Ackr4", "ECs", 3, 1, 2.959658333333334, 8.878975000000001, 0.1465673327253718, 0.1465673327253718, 3, 1, 0.2315733333333333, 0.69472, 0.2371078251520917, 0.2371078251520917, 0.6853779457777778, 6.168401512, 0.03475226150085592, 0.03475226150085592),
    @("M2", "Ccl25", "Ackr4", "FAPs", 3, 1, 2.959658333333334, 8.878975000000001, 0.1465673327253718, 0.1465673327253718, 3, 1, 0.7200953333333334, 2.160286, 0.7373052671097876, 0.7373052671097876, 2.131236154094445, 19.18112538685, 0.1080648664046494, 0.1080648664046494),
    @("M2", "Ccl25", "Ackr4", "sCs", 3, 1, 2.959658333333334, 8.878975000000001, 0.1465673327253718, 0.1465673327253718, 1, 0.3333333333333333, 0.02498966666666666, 0.07496899999999999, 0.02558690773812063, 0.02558690773812063, 0.07396087519722222, 0.665647876775, 0.003750204819866517, 0.003750204819866518),
    @("sCs", "Ccl25", "Ackr4", "ECs", 3, 1, 5.205681333333334, 15.617044, 0.2577942255873873, 0.2577942255873873, 3, 1, 0.2315733333333333, 0.69472, 0.2371078251520917, 0.2371078251520917, 1.205496978631111, 10.84947280768, 0.06112502816579311, 0.06112502816579311),
    @("sCs", "Ccl25", "Ackr4", "FAPs", 3, 1, 5.205681333333334, 15.617044, 0.2577942255873873, 0.2577942255873873, 3, 1, 0.7200953333333334, 2.160286, 0.7373052671097876, 0.7373052671097876, 3.748586834953779, 33.73728151458401, 0.1900730403560694, 0.1900730403560694),
    @("sCs", "Ccl25", "Ackr4", "sCs", 3, 1, 5.205681333333334, 15.617044, 0.2577942255873873, 0.2577942255873873, 1, 0.3333333333333333, 0.02498966666666666, 0.07496899999999999, 0.02558690773812063, 0.02558690773812063, 0.1300882412928889, 1.170794171636, 0.006596157065524734, 0.006596157065524735),
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = $i + 2
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $rowVals.Count; $j++) {
        $colNum = $j + 1
        $ws.Cells.Item($rowNum, $colNum).Value = $rowVals[$j]
    }
}
